$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 51332
$ws.Range("J88").Value = 60598.4
$ws.Range("L88").Value = 60598.4
$ws.Range("N88").Value = -61410.4
$ws.Range("H91").Value = 51332
$ws.Range("J91").Value = 60598.4
$ws.Range("L91").Value = 60598.4
$ws.Range("N91").Value = -63406.4
$ws.Range("H106").Value = 3404
$ws.Range("I106").Value = 3404
$ws.Range("K106").Value = 3404
$ws.Range("M106").Value = -2773
$ws.Range("H113").Value = 6510.4443
$ws.Range("I113").Value = 5120
$ws.Range("K113").Value = 5120
$ws.Range("M113").Value = -1866
$ws.Range("H125").Value = 336665.34
$ws.Range("J125").Value = 336665.34
$ws.Range("L125").Value = 3029988.06
$ws.Range("N125").Value = -3034908.06

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5282.4062
$ws.Range("I32").Value = 1365.56
$ws.Range("J32").Value = 19271.143
$ws.Range("K32").Value = 1365.56
$ws.Range("L32").Value = 19271.143
$ws.Range("M32").Value = -1078.56
$ws.Range("N32").Value = -19845.143
$ws.Range("H122").Value = 1635.6923
$ws.Range("I122").Value = 1553.7778
$ws.Range("J122").Value = 1820
$ws.Range("K122").Value = 4661.3334
$ws.Range("L122").Value = 5460
$ws.Range("M122").Value = -2211.3334
$ws.Range("N122").Value = -10360
$ws.Range("H123").Value = 65000
$ws.Range("J123").Value = 65000
$ws.Range("L123").Value = 65000
$ws.Range("N123").Value = -74800
$ws.Range("H132").Value = 3188.25
$ws.Range("I132").Value = 2721.4707
$ws.Range("K132").Value = 8164.4121
$ws.Range("M132").Value = -5634.4121

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3699
$ws.Range("J20").Value = 3732.1667
$ws.Range("L20").Value = 3732.1667
$ws.Range("N20").Value = -4226.1667
$ws.Range("H105").Value = 3971.2856
$ws.Range("I105").Value = 3971.2856
$ws.Range("K105").Value = 3971.2856
$ws.Range("M105").Value = -2224.2856
$ws.Range("H134").Value = 6061.1665
$ws.Range("J134").Value = 8290.333000000001
$ws.Range("L134").Value = 24870.999
$ws.Range("N134").Value = -29940.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 195.66667
$ws.Range("J7").Value = 500
$ws.Range("L7").Value = 500
$ws.Range("N7").Value = -726
$ws.Range("H22").Value = 966.3333
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 899
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 899
$ws.Range("M22").Value = -650
$ws.Range("N22").Value = -1599
$ws.Range("H31").Value = 1065.9286
$ws.Range("J31").Value = 2011
$ws.Range("L31").Value = 2011
$ws.Range("N31").Value = -2601
$ws.Range("H34").Value = 1065.9286
$ws.Range("J34").Value = 2011
$ws.Range("L34").Value = 2011
$ws.Range("N34").Value = -2415
$ws.Range("H50").Value = 19714.857
$ws.Range("J50").Value = 19426.428
$ws.Range("L50").Value = 19426.428
$ws.Range("N50").Value = -20676.428
$ws.Range("H51").Value = 19999.75
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 19999.75
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 19999.75
$ws.Range("N51").Value = -21471.75
$ws.Range("H60").Value = 23022.75
$ws.Range("H61").Value = 19999.75
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 19999.75
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 19999.75
$ws.Range("N61").Value = -20695.75
$ws.Range("H86").Value = 10888.223
$ws.Range("J86").Value = 11428
$ws.Range("L86").Value = 11428
$ws.Range("N86").Value = -13674
$ws.Range("H89").Value = 10888.223
$ws.Range("J89").Value = 11428
$ws.Range("L89").Value = 57140
$ws.Range("N89").Value = -68372
$ws.Range("H134").Value = 4507.2607
$ws.Range("I134").Value = 4812
$ws.Range("J134").Value = 4493.409
$ws.Range("K134").Value = 14436
$ws.Range("L134").Value = 13480.227
$ws.Range("M134").Value = -11901
$ws.Range("N134").Value = -18550.227
$ws.Range("M51").ClearContents()
$ws.Range("M61").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1175
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 1233.3334
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 3700.0002
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -13780.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1277.8572
$ws.Range("I97").Value = 1277.8572
$ws.Range("K97").Value = 1277.8572
$ws.Range("M97").Value = -781.8571999999999
$ws.Range("H102").Value = 25010.2
$ws.Range("I102").Value = 3734
$ws.Range("K102").Value = 3734
$ws.Range("M102").Value = -2112
$ws.Range("H113").Value = 863.3333
$ws.Range("I113").Value = 795
$ws.Range("K113").Value = 795
$ws.Range("M113").Value = 1375
$ws.Range("H122").Value = 2848.375
$ws.Range("I122").Value = 2257
$ws.Range("J122").Value = 3308.3333
$ws.Range("K122").Value = 6771
$ws.Range("L122").Value = 9924.999899999999
$ws.Range("M122").Value = -4321
$ws.Range("N122").Value = -14824.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12619.647
$ws.Range("I7").Value = 13108.934
$ws.Range("K7").Value = 13108.934
$ws.Range("M7").Value = -12996.934
$ws.Range("H22").Value = 2024.9
$ws.Range("I22").Value = 2124.875
$ws.Range("K22").Value = 2124.875
$ws.Range("M22").Value = -1829.875
$ws.Range("H27").Value = 2024.9
$ws.Range("I27").Value = 2124.875
$ws.Range("K27").Value = 2124.875
$ws.Range("M27").Value = -2017.875
$ws.Range("H126").Value = 12619.647
$ws.Range("I126").Value = 13108.934
$ws.Range("K126").Value = 39326.802
$ws.Range("M126").Value = -36856.802

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1738.5714
$ws.Range("I96").Value = 995
$ws.Range("K96").Value = 995
$ws.Range("M96").Value = 378
$ws.Range("H136").Value = 7197.643
$ws.Range("I136").Value = 7469.857
$ws.Range("J136").Value = 6925.4287
$ws.Range("K136").Value = 22409.571
$ws.Range("L136").Value = 20776.2861
$ws.Range("M136").Value = -19859.571
$ws.Range("N136").Value = -25876.2861
